$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game Skills")

$newRow = 21

$ws.Cells.Item($newRow, 1).Value = "Casting Accuracy"
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = "When using any spell - that does damage - this skill will be used to see if your spell hits or fails. The higher the better chance you have to hit your enemy with your spell. Your casting accuracy skill bonus is used to avoid your spells from being evaded. Casters will use 5% of their focus mod + this skill bonus, other classes will just use the skill bonus."
$ws.Cells.Item($newRow, 4).Value = 999
$ws.Cells.Item($newRow, 10).Value = 1
$ws.Cells.Item($newRow, 11).Value = 1
$ws.Cells.Item($newRow, 12).Value = 0.001
$ws.Cells.Item($newRow, 14).Value = 0
